# "Generate Report for Handoff"
#
# Status moves from "In Translation" to "Ready for handoff" across the
# Overview sheet (zh-cn / de-de columns) and each language sheet's Status
# column, and the corresponding "Latest Handoff/HO Xliff Generate" datetime
# stamps are bumped to the new handoff-generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns + HO Xliff generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-12 14:45:33"

# zh-cn sheet: Status + Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-12 14:45:26"

# de-de sheet: Status + Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-12 14:45:33"

# Widen the status columns so the longer "Ready for handoff" text fits.
$overview.Columns.Item(5).ColumnWidth = 17.22
$overview.Columns.Item(6).ColumnWidth = 17.22
$zhcn.Columns.Item(3).ColumnWidth = 17.22
$dede.Columns.Item(3).ColumnWidth = 17.22
